$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.636417666666667
$ws.Range("H2").Value = 22.909253
$ws.Range("I2").Value = 0.108532481296676
$ws.Range("J2").Value = 0.108532481296676
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 37.01752445727644
$ws.Range("R2").Value = 333.157720115488
$ws.Range("S2").Value = 0.003690833370943947
$ws.Range("T2").Value = 0.003690833370943947

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.636417666666667
$ws.Range("H3").Value = 22.909253
$ws.Range("I3").Value = 0.108532481296676
$ws.Range("J3").Value = 0.108532481296676
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 12.42516427598222
$ws.Range("R3").Value = 111.82647848384
$ws.Range("S3").Value = 0.001238851371657332
$ws.Range("T3").Value = 0.001238851371657333

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.636417666666667
$ws.Range("H4").Value = 22.909253
$ws.Range("I4").Value = 0.108532481296676
$ws.Range("J4").Value = 0.108532481296676
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 1039.092982488626
$ws.Range("R4").Value = 9351.836842397637
$ws.Range("S4").Value = 0.1036027965540747
$ws.Range("T4").Value = 0.1036027965540747

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.103385
$ws.Range("H5").Value = 45.31015499999999
$ws.Range("I5").Value = 0.214656652056136
$ws.Range("J5").Value = 0.214656652056136
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 73.21363864965332
$ws.Range("R5").Value = 658.9227478468799
$ws.Range("S5").Value = 0.007299768007129814
$ws.Range("T5").Value = 0.007299768007129815

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.103385
$ws.Range("H6").Value = 45.31015499999999
$ws.Range("I6").Value = 0.214656652056136
$ws.Range("J6").Value = 0.214656652056136
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 24.57461704426667
$ws.Range("R6").Value = 221.1715533984
$ws.Range("S6").Value = 0.002450212919284463
$ws.Range("T6").Value = 0.002450212919284463

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.103385
$ws.Range("H7").Value = 45.31015499999999
$ws.Range("I7").Value = 0.214656652056136
$ws.Range("J7").Value = 0.214656652056136
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 2055.128733179206
$ws.Range("R7").Value = 18496.15859861286
$ws.Range("S7").Value = 0.2049066711297217
$ws.Range("T7").Value = 0.2049066711297217

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.62086333333334
$ws.Range("H8").Value = 142.86259
$ws.Range("I8").Value = 0.676810866647188
$ws.Range("J8").Value = 0.676810866647188
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.847498666666667
$ws.Range("N8").Value = 14.542496
$ws.Range("O8").Value = 0.03400671694637637
$ws.Range("P8").Value = 0.03400671694637637
$ws.Range("Q8").Value = 230.8420715138489
$ws.Range("R8").Value = 2077.57864362464
$ws.Range("S8").Value = 0.0230161155683026
$ws.Range("T8").Value = 0.0230161155683026

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.62086333333334
$ws.Range("H9").Value = 142.86259
$ws.Range("I9").Value = 0.676810866647188
$ws.Range("J9").Value = 0.676810866647188
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.627093333333334
$ws.Range("N9").Value = 4.88128
$ws.Range("O9").Value = 0.01141456784970118
$ws.Range("P9").Value = 0.01141456784970118
$ws.Range("Q9").Value = 77.48358925724446
$ws.Range("R9").Value = 697.3523033152001
$ws.Range("S9").Value = 0.007725503558759386
$ws.Range("T9").Value = 0.007725503558759386

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 47.62086333333334
$ws.Range("H10").Value = 142.86259
$ws.Range("I10").Value = 0.676810866647188
$ws.Range("J10").Value = 0.676810866647188
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.0707373333333
$ws.Range("N10").Value = 408.212212
$ws.Range("O10").Value = 0.9545787152039225
$ws.Range("P10").Value = 0.9545787152039225
$ws.Range("Q10").Value = 6479.805986216566
$ws.Range("R10").Value = 58318.25387594909
$ws.Range("S10").Value = 0.6460692475201261
$ws.Range("T10").Value = 0.6460692475201261
